# Applies the "Adding labs with a correct ic" edit to the pneumothorax data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values (columns B, C, D for rows 2-12) ---

# Row 2: Cardiac Output
$ws.Range("B2").Value = 5346
$ws.Range("C2").Value = 5309
$ws.Range("D2").Value = 4982

# Row 3: Heart Rate
$ws.Range("C3").Value = 77
$ws.Range("D3").Value = 86

# Row 4: Stroke Volume
$ws.Range("B4").Value = 75
$ws.Range("C4").Value = 69
$ws.Range("D4").Value = 58

# Row 5: Arterial pO2
$ws.Range("C5").Value = 47

# Row 6: Blood Volume
$ws.Range("B6").Value = 5413
$ws.Range("C6").Value = 5413
$ws.Range("D6").Value = 5345

# Row 7: Total Ventilation
$ws.Range("B7").Value = 6.5
$ws.Range("C7").Value = 6.4
$ws.Range("D7").Value = 7.8

# Row 9: Tidal Volume
$ws.Range("B9").Value = 516
$ws.Range("C9").Value = 409
$ws.Range("D9").Value = 462

# Row 12: Respiratory Drive
$ws.Range("B12").Value = 1.1200000000000001
$ws.Range("D12").Value = 2.1800000000000002

# --- Normalize row heights for rows 2-12 to 15.75 ---
for ($r = 2; $r -le 12; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

# --- Update the saved selection/active cell ---
$ws.Range("D6").Select()
